$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new donation history rows (15 & 16) for "denny ariyana", same shape
# as the existing rows: Riwayat ID, User ID, Nama, Nominal, Tanggal Lengkap,
# Donasi ID, Metode Pembayaran.

# Row 15
$ws.Range("A15").Value = 72
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "denny ariyana"
$ws.Range("D15").Value = 300000
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = "transfer"

# Row 16
$ws.Range("A16").Value = 73
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = "denny ariyana"
$ws.Range("D16").Value = 1500000
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = "transfer"

# "Tanggal Lengkap" (E15/E16) needs the literal text "2023-06-12" stored as
# a shared string, matching column E elsewhere in the sheet. Assigning the
# string straight to .Value makes Excel's auto-detection treat it as a date
# serial (and stamps a new number-format style on the cell), so instead we
# produce the text via a throwaway formula cell and bring it over with a
# values-only paste, which preserves it as plain text.
$ws.Range("Z1").Formula = '="2023-06-12"'
$ws.Range("Z1").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
